# Insert a new data row at row 404 (shifting existing rows 404:493 down to 405:494)
# and populate it with a new "Feria Lagunitas de Puerto Montt - Zanahoria" record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a whole row before the current row 404, pushing data down.
$ws.Rows("404:404").Insert()

# Populate the newly inserted row 404 with the new record's values.
$ws.Range("A404").Value = 4
$ws.Range("B404").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C404").Value = "Los Lagos"
$ws.Range("D404").Value = 44943
$ws.Range("E404").Value = 10
$ws.Range("F404").Value = 100114013
$ws.Range("G404").Value = "Zanahoria"
$ws.Range("H404").Value = "Sin especificar"
$ws.Range("I404").Value = "Primera"
$ws.Range("J404").Value = 900
$ws.Range("K404").Value = 13000
$ws.Range("L404").Value = 14000
$ws.Range("M404").Value = 13500
$ws.Range("N404").Value = '$/saco 20 kilos'
$ws.Range("O404").Value = "Chillán"
$ws.Range("P404").Value = 675
$ws.Range("Q404").Value = 20
$ws.Range("R404").Value = "Hortaliza"
